# Apply the "splitting of faglærere into Lektorutdannede / Faglærere / Yrkesfaglærere"
# edit, and restore the PPU / PPU Yrkesfag blocks (pushed further down the sheet)
# with their original, unmodified figures.
#
# Layout before the edit (sheet "Sheet1"):
#   Row 1        : header
#   Rows 2-22    : Barnehagelærere   (merged A2:A22)
#   Rows 23-43   : Grunnskolelærere  (merged A23:A43)
#   Rows 44-64   : Faglærere         (merged A44:A64)
#   Rows 65-85   : PPU               (merged A65:A85)
#   Rows 86-106  : PPU Yrkesfag      (merged A86:A106)
#
# Layout after the edit:
#   Rows 44-64   : Lektorutdannede   (merged A44:A64)   <- NEW figures
#   Rows 65-85   : Faglærere         (merged A65:A85)   <- NEW figures (renamed/recomputed)
#   Rows 86-106  : Yrkesfaglærere    (merged A86:A106)  <- NEW figures (renamed/recomputed)
#   Rows 107-127 : PPU               (merged A107:A127) <- original PPU figures, shifted down
#   Rows 128-148 : PPU Yrkesfag      (merged A128:A148) <- original PPU Yrkesfag figures, shifted down

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @("Lektorutdannede", 2020, 5531, 5531, 0),
    @($null, 2021, 5860, 5538, 322),
    @($null, 2022, 6218, 5558, 660),
    @($null, 2023, 6585, 5572, 1013),
    @($null, 2024, 6955, 5580, 1374),
    @($null, 2025, 7324, 5590, 1734),
    @($null, 2026, 7687, 5597, 2090),
    @($null, 2027, 8047, 5597, 2449),
    @($null, 2028, 8396, 5579, 2817),
    @($null, 2029, 8740, 5554, 3186),
    @($null, 2030, 9072, 5528, 3544),
    @($null, 2031, 9395, 5509, 3886),
    @($null, 2032, 9709, 5496, 4213),
    @($null, 2033, 10015, 5491, 4524),
    @($null, 2034, 10319, 5478, 4842),
    @($null, 2035, 10620, 5460, 5160),
    @($null, 2036, 10914, 5447, 5467),
    @($null, 2037, 11206, 5448, 5758),
    @($null, 2038, 11498, 5459, 6039),
    @($null, 2039, 11790, 5471, 6319),
    @($null, 2040, 12076, 5486, 6590),
    @("Faglærere", 2020, 5534, 5534, 0),
    @($null, 2021, 5860, 5541, 319),
    @($null, 2022, 6218, 5561, 657),
    @($null, 2023, 6585, 5575, 1010),
    @($null, 2024, 6955, 5583, 1372),
    @($null, 2025, 7324, 5593, 1731),
    @($null, 2026, 7687, 5600, 2087),
    @($null, 2027, 8047, 5600, 2447),
    @($null, 2028, 8396, 5582, 2814),
    @($null, 2029, 8740, 5557, 3183),
    @($null, 2030, 9072, 5531, 3541),
    @($null, 2031, 9395, 5512, 3883),
    @($null, 2032, 9709, 5499, 4210),
    @($null, 2033, 10015, 5494, 4521),
    @($null, 2034, 10319, 5481, 4839),
    @($null, 2035, 10620, 5463, 5157),
    @($null, 2036, 10914, 5450, 5464),
    @($null, 2037, 11206, 5451, 5755),
    @($null, 2038, 11498, 5461, 6037),
    @($null, 2039, 11790, 5474, 6316),
    @($null, 2040, 12076, 5489, 6587),
    @("Yrkesfaglærere", 2020, 5535, 5535, 0),
    @($null, 2021, 5860, 5542, 318),
    @($null, 2022, 6218, 5562, 656),
    @($null, 2023, 6585, 5576, 1010),
    @($null, 2024, 6955, 5584, 1371),
    @($null, 2025, 7324, 5594, 1730),
    @($null, 2026, 7687, 5601, 2086),
    @($null, 2027, 8047, 5601, 2446),
    @($null, 2028, 8396, 5583, 2813),
    @($null, 2029, 8740, 5557, 3182),
    @($null, 2030, 9072, 5532, 3540),
    @($null, 2031, 9395, 5513, 3882),
    @($null, 2032, 9709, 5500, 4209),
    @($null, 2033, 10015, 5495, 4520),
    @($null, 2034, 10319, 5482, 4838),
    @($null, 2035, 10620, 5464, 5156),
    @($null, 2036, 10914, 5451, 5463),
    @($null, 2037, 11206, 5452, 5754),
    @($null, 2038, 11498, 5462, 6036),
    @($null, 2039, 11790, 5475, 6315),
    @($null, 2040, 12076, 5490, 6586),
    @("PPU", 2020, 30789, 30789, 0),
    @($null, 2021, 30793, 30824, -31),
    @($null, 2022, 30890, 30930, -40),
    @($null, 2023, 31004, 31003, 1),
    @($null, 2024, 31130, 31043, 87),
    @($null, 2025, 31256, 31097, 159),
    @($null, 2026, 31384, 31132, 252),
    @($null, 2027, 31496, 31124, 373),
    @($null, 2028, 31594, 31006, 588),
    @($null, 2029, 31672, 30845, 827),
    @($null, 2030, 31717, 30679, 1038),
    @($null, 2031, 31726, 30555, 1170),
    @($null, 2032, 31699, 30462, 1238),
    @($null, 2033, 31633, 30416, 1217),
    @($null, 2034, 31553, 30314, 1238),
    @($null, 2035, 31437, 30184, 1254),
    @($null, 2036, 31322, 30078, 1244),
    @($null, 2037, 31190, 30063, 1127),
    @($null, 2038, 31069, 30099, 970),
    @($null, 2039, 30946, 30146, 800),
    @($null, 2040, 30834, 30213, 621),
    @("PPU Yrkesfag", 2020, 13299, 13299, 0),
    @($null, 2021, 13226, 13325, -100),
    @($null, 2022, 13190, 13386, -196),
    @($null, 2023, 13140, 13443, -304),
    @($null, 2024, 13091, 13490, -399),
    @($null, 2025, 13024, 13559, -535),
    @($null, 2026, 12945, 13624, -679),
    @($null, 2027, 12866, 13657, -792),
    @($null, 2028, 12763, 13620, -857),
    @($null, 2029, 12664, 13556, -893),
    @($null, 2030, 12555, 13480, -924),
    @($null, 2031, 12434, 13426, -992),
    @($null, 2032, 12321, 13387, -1066),
    @($null, 2033, 12191, 13373, -1182),
    @($null, 2034, 12080, 13314, -1234),
    @($null, 2035, 11961, 13229, -1267),
    @($null, 2036, 11852, 13153, -1301),
    @($null, 2037, 11748, 13129, -1381),
    @($null, 2038, 11635, 13132, -1498),
    @($null, 2039, 11542, 13138, -1597),
    @($null, 2040, 11443, 13155, -1712)
)

$startRow = 44
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $entry = $rowsData[$i]
    $r = $startRow + $i

    if ($entry[0] -ne $null) {
        $ws.Cells.Item($r, 1).Value = $entry[0]
    }
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
    $ws.Cells.Item($r, 5).Value = $entry[4]
}

# Re-create the vertical "group label" merges for the two blocks that moved
# down (do this BEFORE styling, so every row still gets its own full-box
# border below instead of Excel collapsing the merged column to one outline).
$ws.Range("A107:A127").Merge()
$ws.Range("A128:A148").Merge()

# Column A / B in this table use the bold, thin-bordered, centered/top-aligned
# style (style index 1 in the original workbook) on EVERY row (each cell gets
# its own full box border, not just the outline of the whole block). Apply the
# same look to the two newly-created blocks (rows 107-148) so they match the
# rest of the table. Borders are set per-row because Range.Borders only draws
# the outline of a multi-row range, not a box around every row in it.
for ($r = 107; $r -le 148; $r++) {
    $rowRange = $ws.Range("A" + $r + ":B" + $r)
    $rowRange.Font.Bold = $true
    $rowRange.HorizontalAlignment = -4108
    $rowRange.VerticalAlignment = -4160
    $rowRange.Borders.LineStyle = 1
    $rowRange.Borders.Weight = 2
}

Write-Output "done"
